$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B column label (rows 2-10) from "Diferença 2023/04 - 2022/04" to "Diferença 2024/01 - 2023/01"
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 2).Value = "Diferença 2024/01 - 2023/01"
}

# Row 3: Piauí -> Rondônia
$ws.Range("A3").Value = "Rondônia"
# Row 4: Ceará -> Rio Grande do Sul
$ws.Range("A4").Value = "Rio Grande do Sul"
# Row 5: Amapá -> Mato Grosso do Sul
$ws.Range("A5").Value = "Mato Grosso do Sul"
# Row 6: Mato Grosso do Sul -> Santa Catarina
$ws.Range("A6").Value = "Santa Catarina"
# Row 7: Rondônia -> Bahia
$ws.Range("A7").Value = "Bahia"

# Update C column values
$ws.Range("C2").Value = 0.7999999999999998
$ws.Range("C3").Value = 0.5
$ws.Range("C4").Value = 0.3999999999999995
$ws.Range("C5").Value = 0.2000000000000002
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = -0.4000000000000004
$ws.Range("C8").Value = -1.9
$ws.Range("C9").Value = -1.1
$ws.Range("C10").Value = -0.9000000000000004

# Update D8: 20º -> 27º
$ws.Range("D8").Value = "27º"
